$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (40 and 41) following the existing pattern:
# column A: date serial 45920 (2025-09-20), formatted like existing date cells
# column B: station name (shared text), matching rows for "四方坪站" and "高岭站"

$ws.Cells.Item(39, 1).Copy()
$ws.Cells.Item(40, 1).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(41, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(40, 1).Value = 45920
$ws.Cells.Item(40, 2).Value = "四方坪站"

$ws.Cells.Item(41, 1).Value = 45920
$ws.Cells.Item(41, 2).Value = "高岭站"

# Update selection to match the new state (C42)
$ws.Range("C42").Select()
